# Applies the edits described by the diff to usa_mls_2023 workbook:
#  1) Row 169 and Row 173: the match data in columns F:V is swapped
#     (row169 <-> row173), row index/meta columns A:E untouched.
#  2) Rows 412, 413, 414: the match data in columns F:V is rotated
#     (412<-413, 413<-414, 414<-412), row index/meta columns A:E untouched.
#  3) Two brand-new rows (434 and 435) are appended at the bottom of the
#     sheet, extending the used range from A1:V433 to A1:V435.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# --- 1) swap F:V between row 169 and row 173 -------------------------------
$row169 = Get-RowValues 169
$row173 = Get-RowValues 173

Set-RowValues 169 $row173
Set-RowValues 173 $row169

# --- 2) rotate F:V across rows 412, 413, 414 --------------------------------
# after: row412 <- old row413, row413 <- old row414, row414 <- old row412
$row412 = Get-RowValues 412
$row413 = Get-RowValues 413
$row414 = Get-RowValues 414

Set-RowValues 412 $row413
Set-RowValues 413 $row414
Set-RowValues 414 $row412

# --- 3) append two new rows (434 and 435) -----------------------------------
# Copy formatting (styles) from the last existing data row (433) so the
# new rows match the sheet's look (bold/border/centered index column,
# date-formatted match-date column, etc.)
$ws.Range("A433:V433").Copy()
$ws.Range("A434:V435").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 434
$ws.Range("A434").Value = 433
$ws.Range("B434").Value = "usa"
$ws.Range("C434").Value = "mls"
# D434 ("2023") must be stored as TEXT (like every other row), not a number.
# Excel auto-detects an all-digit string as numeric, so force a Text format
# first, then strip the format back to the sheet's normal (unstyled) look
# by pasting formats from a neighbouring plain-text cell.
$ws.Range("D434").NumberFormat = "@"
$ws.Range("D434").Value = "2023"
$ws.Range("E434").Value = 45194.0625
$ws.Range("F434").Value = "Orlando City"
$ws.Range("G434").Value = 1
$ws.Range("H434").Value = "Inter Miami"
$ws.Range("I434").Value = 1
$ws.Range("J434").Value = 1.96
$ws.Range("K434").Value = "21/09/2023 08:48"
$ws.Range("L434").Value = 1.72
$ws.Range("M434").Value = "25/09/2023 01:28"
$ws.Range("N434").Value = 3.76
$ws.Range("O434").Value = "21/09/2023 08:48"
$ws.Range("P434").Value = 4.43
$ws.Range("Q434").Value = "25/09/2023 01:29"
$ws.Range("R434").Value = 3.78
$ws.Range("S434").Value = "21/09/2023 08:48"
$ws.Range("T434").Value = 4.44
$ws.Range("U434").Value = "25/09/2023 01:29"
$ws.Range("V434").Value = "https://www.betexplorer.com/football/usa/mls/orlando-city-inter-miami/Yw0V68Ii/"

# Row 435
$ws.Range("A435").Value = 434
$ws.Range("B435").Value = "usa"
$ws.Range("C435").Value = "mls"
$ws.Range("D435").NumberFormat = "@"
$ws.Range("D435").Value = "2023"
$ws.Range("E435").Value = 45194.14583333334
$ws.Range("F435").Value = "Austin FC"
$ws.Range("G435").Value = 3
$ws.Range("H435").Value = "Los Angeles Galaxy"
$ws.Range("I435").Value = 3
$ws.Range("J435").Value = 1.88
$ws.Range("K435").Value = "21/09/2023 03:43"
$ws.Range("L435").Value = 2
$ws.Range("M435").Value = "25/09/2023 02:59"
$ws.Range("N435").Value = 4.12
$ws.Range("O435").Value = "21/09/2023 03:43"
$ws.Range("P435").Value = 4.24
$ws.Range("Q435").Value = "25/09/2023 02:59"
$ws.Range("R435").Value = 3.75
$ws.Range("S435").Value = "21/09/2023 03:43"
$ws.Range("T435").Value = 3.42
$ws.Range("U435").Value = "25/09/2023 02:54"
$ws.Range("V435").Value = "https://www.betexplorer.com/football/usa/mls/austin-fc-los-angeles-galaxy/G2kz6SXc/"

# Now that the "2023" text values are safely stored, restore D434/D435 to the
# sheet's normal (unstyled / General) look, matching every other row's D cell.
$ws.Range("B434:B435").Copy()
$ws.Range("D434:D435").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Done."
